$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-35 down to 28-36,
# and fill it with this week's new record for "Agrícola del Norte S.A. de Arica".
$ws.Rows.Item(27).Insert()

$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44825
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112045
$ws.Cells.Item(27, 7).Value = "Zapallo"
$ws.Cells.Item(27, 8).Value = "Camote"
$ws.Cells.Item(27, 9).Value = "1a nueva(o)"
$ws.Cells.Item(27, 10).Value = 1200
$ws.Cells.Item(27, 11).Value = 900
$ws.Cells.Item(27, 12).Value = 930
$ws.Cells.Item(27, 13).Value = 915
$ws.Cells.Item(27, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(27, 15).Value = "Perú"
$ws.Cells.Item(27, 16).Value = 915
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = "Hortaliza"
